$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (V42_TEAGER_ENSEMBLE)
$ws.Range("A40").Value = "V42_TEAGER_ENSEMBLE"
$ws.Range("B40").Value = 339
$ws.Range("C40").Value = 0

# Copy formatting (style) from the row above (A36, which carries the bold/custom style)
$ws.Range("A36").Copy()
$ws.Range("A40").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the selection to the newly added cell (matches the saved sheetView state)
$ws.Range("C40").Select()
